$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "Datos actualizados" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 18:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 929841
$ws.Range("C4").Value = 4609
$ws.Range("E4").Value = 766494
$ws.Range("G4").Value = 650
$ws.Range("H4").Value = 52843

# Row 6 - Italia
$ws.Range("B6").Value = 195351
$ws.Range("C6").Value = 2357
$ws.Range("D6").Value = 63120
$ws.Range("E6").Value = 105847
$ws.Range("F6").Value = 2102
$ws.Range("G6").Value = 415
$ws.Range("H6").Value = 26384

# Row 33 - Polonia
$ws.Range("B33").Value = 11273
$ws.Range("C33").Value = 381
$ws.Range("E33").Value = 8623
$ws.Range("G33").Value = 30
$ws.Range("H33").Value = 524

# Row 55 - Marruecos
$ws.Range("B55").Value = 3897
$ws.Range("C55").Value = 139
$ws.Range("D55").Value = 537
$ws.Range("E55").Value = 3201

# Row 108 - Sri Lanka
$ws.Range("B108").Value = 452
$ws.Range("C108").Value = 35
$ws.Range("E108").Value = 327
